# Weekly data refresh for "Fruta, Macroferia Regional de Talca - Chirimoya":
# two new price records (2021-10-25) are inserted at the top of the data block,
# pushing the existing rows down by two and appending the two oldest records
# (previously rows 43-44) to the bottom as rows 45-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row {
    param($r, $d, $l, $m, $n, $o, $p, $s)
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 19).Value = $s
}

function Set-NewRowConstants {
    param($r)
    $ws.Cells.Item($r, 1).Value = 5
    $ws.Cells.Item($r, 2).Value = "Macroferia Regional de Talca"
    $ws.Cells.Item($r, 3).Value = "Maule"
    $ws.Cells.Item($r, 5).Value = 7
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100107
    $ws.Cells.Item($r, 8).Value = "Otros"
    $ws.Cells.Item($r, 9).Value = 100107002
    $ws.Cells.Item($r, 10).Value = "Chirimoya"
    $ws.Cells.Item($r, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"
    $ws.Cells.Item($r, 18).Value = "Provincia de Limarí"
    $ws.Cells.Item($r, 20).Value = 10
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}

Set-Row 18 44494 "Primera" 150 25000 25000 25000 2500
Set-Row 19 44494 "Segunda" 50 23000 23000 23000 2300
Set-Row 20 44484 "Primera" 120 25000 25000 25000 2500
Set-Row 21 44484 "Segunda" 100 22000 22000 22000 2200
Set-Row 22 44447 "Especial" 50 32000 32000 32000 3200
Set-Row 23 44421 "Especial" 30 35000 35000 35000 3500
Set-Row 24 44454 "Especial" 320 30000 30000 30000 3000
Set-Row 25 44454 "Primera" 300 28000 28000 28000 2800
Set-Row 26 44467 "Especial" 100 30000 30000 30000 3000
Set-Row 27 44467 "Primera" 100 28000 28000 28000 2800
Set-Row 28 44473 "Primera" 200 28000 28000 28000 2800
Set-Row 29 44434 "Especial" 60 30000 30000 30000 3000
Set-Row 30 44445 "Primera" 250 28000 30000 29200 2920
Set-Row 31 44489 "Especial" 50 27000 27000 27000 2700
Set-Row 32 44489 "Primera" 50 25000 25000 25000 2500
Set-Row 33 44475 "Primera" 200 28000 28000 28000 2800
Set-Row 34 44474 "Especial" 150 30000 30000 30000 3000
Set-Row 35 44448 "Especial" 100 30000 30000 30000 3000
Set-Row 36 44448 "Primera" 80 28000 28000 28000 2800
Set-Row 37 44453 "Especial" 135 30000 30000 30000 3000
Set-Row 38 44435 "Especial" 160 30000 30000 30000 3000
Set-Row 39 44476 "Especial" 100 30000 30000 30000 3000
Set-Row 40 44468 "Especial" 250 30000 30000 30000 3000
Set-Row 41 44491 "Primera" 100 25000 25000 25000 2500
Set-Row 42 44466 "Especial" 110 30000 30000 30000 3000
Set-Row 43 44438 "Primera" 100 30000 30000 30000 3000
Set-Row 44 44461 "Especial" 150 30000 30000 30000 3000
Set-NewRowConstants 45
Set-Row 45 44461 "Primera" 100 25000 25000 25000 2500
Set-NewRowConstants 46
Set-Row 46 44432 "Especial" 70 30000 30000 30000 3000
